$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BOM had a single "PIN HEADER" row (SV1/ICSP/MA03-2) whose part-number
# cell listed both the male and female Digikey part numbers together.
# Split it into two rows: one for the female header, one for the male header.

# Remove the existing hyperlink on F22 (https://www.sparkfun.com/products/8665)
# before the row insert shifts things around, then re-add it afterwards at
# its new location.
[void]$ws.Range("F22").Hyperlinks.Delete()

# Insert a new row above the old row 18 (SV2,SV3 / ICSP...) so the split
# "PIN HEADER" entry gets its own second row; everything below moves down
# by one row.
[void]$ws.Rows.Item(18).Insert()

# Row 17: now the "female" pin header entry (drop the unit cost - the price
# moves to the new male row instead).
$ws.Range("E17").Value = "PIN HEADER female"
$ws.Range("F17").Value = "3M9525-ND "
[void]$ws.Range("G17").ClearContents()

# Row 18 (new): the "male" pin header entry, carrying the price that used
# to live on the combined row.
$ws.Range("A18").Value = "SV1"
$ws.Range("B18").Value = "ICSP"
$ws.Range("C18").Value = "MA03-2"
$ws.Range("D18").Value = "MA03-2"
$ws.Range("E18").Value = "PIN HEADER male"
$ws.Range("F18").Value = "3M9459-ND "
$ws.Range("G18").Value = 0.42

# Re-create the hyperlink, now on F23 (it was F22 before the new row).
[void]$ws.Hyperlinks.Add($ws.Range("F23"), "https://www.sparkfun.com/products/8665")

# The named range that captured the BOM body grows by one row.
$wb.Names.Item("ControlModule_BOM").RefersTo = "=Sheet1!`$A`$1:`$E`$23"

# Restore a sensible selection matching the saved file (whole BOM table
# selected, active cell at its bottom-right corner).
[void]$ws.Range("A4:G23").Select()
